$d = $word.ActiveDocument
$xmlHeader = '<?xml version="1.0" encoding="UTF-8"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Edit 21: paragraphs 111..111
$p1 = $d.Paragraphs(111)
$p2 = $d.Paragraphs(111)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag21 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Unknown </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>value scan</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag21 + $xmlFooter)

# Edit 20: paragraphs 105..105
$p1 = $d.Paragraphs(105)
$p2 = $d.Paragraphs(105)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag20 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Memory protection filters (read, write, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>exec</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>).</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag20 + $xmlFooter)

# Edit 19: paragraphs 94..95
$p1 = $d.Paragraphs(94)
$p2 = $d.Paragraphs(95)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag19 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>NOP/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>UnNOP</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> support.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>ASM searching API with ‘wildcards’. (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MetaASM</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>?)</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag19 + $xmlFooter)

# Edit 18: paragraphs 90..91
$p1 = $d.Paragraphs(90)
$p2 = $d.Paragraphs(91)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag18 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Helper functions such as </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>FindExport</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>FindImport</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>HasDataDir</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GetArchitecture</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"></w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IsDotNet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"></w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GetPDB</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> etc.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Test against pathological cases such as </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Corkami</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tests.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag18 + $xmlFooter)

# Edit 17: paragraphs 83..83
$p1 = $d.Paragraphs(83)
$p2 = $d.Paragraphs(83)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag17 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Load </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> directory.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag17 + $xmlFooter)

# Edit 16: paragraphs 74..75
$p1 = $d.Paragraphs(74)
$p2 = $d.Paragraphs(75)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag16 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">E.g. </w:t></w:r><w:r><w:t xml:space="preserve">Check </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NumberOfRvaAndSizes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>NtHeaders</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> before attempting to retrieve a data dir.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Cache base pointers </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> rather than retrieving it manually in every getter/setter. Slightly less ‘robust’, but due to the typically ‘read-only’ nature of the data this is the expected behaviour in all known cases anyway.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag16 + $xmlFooter)

# Edit 15: paragraphs 72..72
$p1 = $d.Paragraphs(72)
$p2 = $d.Paragraphs(72)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag15 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Note: May cause problems when copying ‘</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PeFile</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>’ type.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag15 + $xmlFooter)

# Edit 14: paragraphs 69..70
$p1 = $d.Paragraphs(69)
$p2 = $d.Paragraphs(70)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag14 = @'
<w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>PeLib</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Investigate use of virtual functions for file </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> memory</w:t></w:r><w:r><w:t xml:space="preserve"> access</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RvaToVa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag14 + $xmlFooter)

# Edit 13: paragraphs 64..65
$p1 = $d.Paragraphs(64)
$p2 = $d.Paragraphs(65)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag13 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Detect cases where hooking may overflow past the end of a function, and fail. (Provide policy or flag to allow overriding this behaviour.) Examples may be instructions such as ‘</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 3’, ‘ret’, ‘</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>jmp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>’, etc.</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>FindPattern</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$rng.InsertXML($xmlHeader + $frag13 + $xmlFooter)

# Edit 12: paragraphs 59..59
$p1 = $d.Paragraphs(59)
$p2 = $d.Paragraphs(59)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag12 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Uncopyable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, so make moveable.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag12 + $xmlFooter)

# Edit 11: paragraphs 56..56
$p1 = $d.Paragraphs(56)
$p2 = $d.Paragraphs(56)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag11 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Improved relative instruction rebuilding (including conditionals).</w:t></w:r><w:r><w:t xml:space="preserve"></w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>x64</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> has far more IP relative instructions than x86.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag11 + $xmlFooter)

# Edit 10: paragraphs 53..53
$p1 = $d.Paragraphs(53)
$p2 = $d.Paragraphs(53)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag10 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Class function hooking (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ecx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> preservation).</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag10 + $xmlFooter)

# Edit 9: paragraphs 50..51
$p1 = $d.Paragraphs(50)
$p2 = $d.Paragraphs(51)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag9 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Get address of Kernel32</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>!LoadLibrary</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> ‘manually’ rather than using a local </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GetProcAddress</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:t>pointer arithmetic.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Whilst this works in all normal cases, it doesn’t work when the target has shims enabled which hook </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LoadLibrary</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag9 + $xmlFooter)

# Edit 8: paragraphs 41..42
$p1 = $d.Paragraphs(41)
$p2 = $d.Paragraphs(42)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag8 = @'
<w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>ManualMap</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Exception handling support under x86 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SafeSEH</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and x64.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag8 + $xmlFooter)

# Edit 7: paragraphs 32..35
$p1 = $d.Paragraphs(32)
$p2 = $d.Paragraphs(35)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag7 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Remote code ‘emulator’.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Custom </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GetModuleHandle</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GetProcAddress</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, etc.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="2160" w:hanging="2160"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>MemoryMgr</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">‘Unchecked’ read/write </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> functions designed for speed and use in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReadString</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"></w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> where you only want to check page protections once, then forget about it.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag7 + $xmlFooter)

# Edit 6: paragraphs 29..30
$p1 = $d.Paragraphs(29)
$p2 = $d.Paragraphs(30)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag6 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Helper service to run </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>HadesMem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tools as ‘SYSTEM’ for when </w:t></w:r><w:r><w:t>manipulating</w:t></w:r><w:r><w:t xml:space="preserve"> certain protected/critical processes (running in separate desktops, sessions, etc.).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Debugging</w:t></w:r><w:r><w:t xml:space="preserve"> APIs (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"></w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>hw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"></w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, conditional </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>bp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, single stepping, stack trace, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag6 + $xmlFooter)

# Edit 5: paragraphs 27..27
$p1 = $d.Paragraphs(27)
$p2 = $d.Paragraphs(27)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag5 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Important! Ensure -</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fno</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">-strict-aliasing is used under GCC as it seems </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Boost.Python</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> has aliasing violations which cause spurious </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>segfaults</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and other issues.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag5 + $xmlFooter)

# Edit 4: paragraphs 23..24
$p1 = $d.Paragraphs(23)
$p2 = $d.Paragraphs(24)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag4 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Merge headers where appropriate (e.g. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>module_list</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>module_iterator</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Check whether Read/Write APIs and other templates should be using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>std</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>::</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>remove_cv</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> anywhere in the type detection/transformation.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Rewrite </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ReadVector</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>WriteVector</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to use </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>true_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>false_type</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> overloading rather than enable_if.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$rng.InsertXML($xmlHeader + $frag4 + $xmlFooter)

# Edit 3: paragraphs 17..17
$p1 = $d.Paragraphs(17)
$p2 = $d.Paragraphs(17)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag3 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Extra sanity checks to ensure not only that functions return without exception, but also that the returned data is valid. (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>e.g</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GetName</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ImpThunk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">.) </w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag3 + $xmlFooter)

# Edit 2: paragraphs 10..11
$p1 = $d.Paragraphs(10)
$p2 = $d.Paragraphs(11)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag2 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Document preconditions and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>postconditions</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Better </w:t></w:r><w:r><w:t>annotations (</w:t></w:r><w:r><w:t xml:space="preserve">warnings, notes, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">). </w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag2 + $xmlFooter)

# Edit 1: paragraphs 4..4
$p1 = $d.Paragraphs(4)
$p2 = $d.Paragraphs(4)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag1 = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Todo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> list</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
'@
$rng.InsertXML($xmlHeader + $frag1 + $xmlFooter)

# Edit 0: paragraphs 1..1
$p1 = $d.Paragraphs(1)
$p2 = $d.Paragraphs(1)
$rng = $d.Range($p1.Range.Start, $p2.Range.End)
$frag0 = @'
<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Todo</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$rng.InsertXML($xmlHeader + $frag0 + $xmlFooter)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"